$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# This sheet stores numeric-looking prices as plain text (e.g. "1.00", "0.370")
# rather than real numbers, so that trailing/insignificant zeros are preserved.
# Whenever the new price text would otherwise be auto-recognized by Excel as a
# number, force the cell to Text format first and clear the format again right
# after the write so the cell's style index is left untouched (matches source).

# Row 39/40: Kaspa and WhiteBITCoin swap places (rank index in column A stays the same)
$ws.Range("B39").Value = "Kaspa"
$ws.Range("C39").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.130"
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = "  +2.02%  "

$ws.Range("B40").Value = "WhiteBITCoin"
$ws.Range("C40").Value = "https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "22.23"
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = "  -0.85%  "

# Price / Volume updates across all rows
$ws.Range("D2").Value = "88.064.08"
$ws.Range("E2").Value = "  -3.13%  "
$ws.Range("D3").Value = "3.057.73"
$ws.Range("E3").Value = "  -4.76%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.00"
$ws.Range("D4").ClearFormats()
$ws.Range("E4").Value = "  -0.10%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "210.67"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  -4.11%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "618.09"
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = "  -4.73%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.370"
$ws.Range("D7").ClearFormats()
$ws.Range("E7").Value = "  -7.38%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.795"
$ws.Range("D8").ClearFormats()
$ws.Range("E8").Value = "  +12.48%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "1.00"
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value = "  -0.02%  "
$ws.Range("D10").Value = "3.053.46"
$ws.Range("E10").Value = "  -4.85%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.593"
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = "  +2.44%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.180"
$ws.Range("D12").ClearFormats()
$ws.Range("E12").Value = "  -1.30%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.0000238"
$ws.Range("D13").ClearFormats()
$ws.Range("E13").Value = "  -8.12%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "5.26"
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = "  -3.21%  "
$ws.Range("D15").Value = "87.782.89"
$ws.Range("E15").Value = "  -3.23%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "31.96"
$ws.Range("D16").ClearFormats()
$ws.Range("E16").Value = "  -4.58%  "
$ws.Range("D17").Value = "3.613.36"
$ws.Range("E17").Value = "  -4.93%  "
$ws.Range("D18").Value = "3.063.50"
$ws.Range("E18").Value = "  -5.47%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "3.31"
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = "  -2.13%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.0000204"
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = "  -9.32%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "13.30"
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = "  -1.73%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "419.41"
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = "  -5.34%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "8.15"
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = "  -6.35%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "4.90"
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = "  -3.96%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "5.47"
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = "  +3.21%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "11.81"
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = "  -0.67%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "81.62"
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = "  -0.48%  "
$ws.Range("D28").Value = "3.195.96"
$ws.Range("E28").Value = "  -5.55%  "
$ws.Range("E29").Value = "  +0.09%  "
$ws.Range("E30").Value = "  +7.98%  "
$ws.Range("E31").Value = "  +8.07%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "8.03"
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = "  -5.44%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "505.98"
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = "  -7.23%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "3.62"
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = "  -13.54%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "6.70"
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = "  -6.25%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.24"
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = "  -4.93%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.79"
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = "  -7.91%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "22.10"
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = "  -2.27%  "
$ws.Range("E41").Value = "  +0.09%  "
$ws.Range("E42").Value = "  +0.01%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.360"
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = "  -4.19%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "147.58"
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = "  +0.69%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "1.80"
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = "  -7.87%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.131"
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = "  +3.68%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "43.28"
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = "  -4.22%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.0675"
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = "  +10.49%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.20"
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = "  -3.73%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "157.21"
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = "  -9.63%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.698"
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = "  -8.85%  "

Write-Output "Applied cryptos update"
